$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for "Approach"
$ws.Range("L2").Value = "Approach"

# New row 8 data: Majority element (duplicate detection / array topic)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Majority element"
$ws.Range("D8").Value = "easy"
$ws.Range("F8").Value = "Yes"
$ws.Range("H8").Value = "Array"
$ws.Range("J8").Value = "NO"
$ws.Range("K8").Value = "Yes"
$ws.Range("L8").Value = "Moore's voting algo."

# Add hyperlink on C8 similar to other question links in the sheet
$ws.Hyperlinks.Add($ws.Range("C8"), "https://leetcode.com/problems/majority-element/description/")
$ws.Range("C8").Style = "Hyperlink"

# Update selection to mirror the authored state
$ws.Range("L10").Select()
